$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New tracker rows appended at the end (rows 62 and 63)
$ws.Range("A62").Value = "G1"
$ws.Range("B62").Value = "Test1"
$ws.Range("C62").Value = 45891
$ws.Range("D62").Value = 0.7493421469649953
$ws.Range("E62").Value = 0
$ws.Range("F62").Value = -0.01

$ws.Range("A63").Value = "G2"
$ws.Range("B63").Value = "sedrftgyhuioygtfrd"
$ws.Range("C63").Value = 45891
$ws.Range("D63").Value = 0.7493421469649953
$ws.Range("E63").Value = 0
$ws.Range("F63").Value = -0.01

# Match the date formatting style used by the other Date column cells
$ws.Range("C62").NumberFormat = $ws.Range("C61").NumberFormat
$ws.Range("C63").NumberFormat = $ws.Range("C61").NumberFormat
